$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PVTStL")

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

$ws.Activate()
$ws.Range("D6").Select()
$wb.Worksheets.Item("About").Activate()
